{"js": "// Apply three edits described by the diff:\n// 1. Append a new run \"DONE\" to the \"Show example of phishing email \" paragraph.\n// 2. Append a new run \" \" (space) to the \"... users name\" paragraph.\n// 3. After the \"Ask Rory about simulating PayPal Login\" paragraph, insert a new\n//    empty paragraph followed by a new paragraph with text\n//    \"Add info button to reappear instructions \".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet showExamplePara = null;\nlet usersNamePara = null;\nlet askRoryPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Show example of phishing email\") !== -1) {\n    showExamplePara = p;\n  } else if (/\\busers name\\s*$/.test(t)) {\n    usersNamePara = p;\n  } else if (t.indexOf(\"Ask Rory about simulating PayPal Login\") !== -1) {\n    askRoryPara = p;\n  }\n}\n\nif (!showExamplePara) {\n  throw new Error(\"Could not find 'Show example of phishing email' paragraph\");\n}\nif (!usersNamePara) {\n  throw new Error(\"Could not find '... users name' paragraph\");\n}\nif (!askRoryPara) {\n  throw new Error(\"Could not find 'Ask Rory about simulating PayPal Login' paragraph\");\n}\n\n// 1. Add a \"DONE\" run at the end of the \"Show example of phishing email\" paragraph.\nshowExamplePara.insertText(\"DONE\", Word.InsertLocation.end);\n\n// 2. Add a trailing space run at the end of the \"... users name\" paragraph.\nusersNamePara.insertText(\" \", Word.InsertLocation.end);\n\n// 3. Insert a blank paragraph and a new task paragraph after the \"Ask Rory...\" paragraph.\nconst infoPara = askRoryPara.insertParagraph(\n  \"Add info button to reappear instructions \",\n  Word.InsertLocation.after\n);\ninfoPara.insertParagraph(\"\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Apply three edits described by the diff:\n# 1. Append a new run \"DONE\" to the \"Show example of phishing email \" paragraph.\n# 2. Append a new run \" \" (space) to the \"... users name\" paragraph.\n# 3. After the \"Ask Rory about simulating PayPal Login\" paragraph, insert a new\n#    empty paragraph followed by a new paragraph with text\n#    \"Add info button to reappear instructions \".\n\n$d = $word.ActiveDocument\n\n# 1. Add a \"DONE\" run at the end of the \"Show example of phishing email\" paragraph.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Show example of phishing email*\") {\n        $r = $p.Range\n        $r.End = $r.End - 1   # exclude the paragraph mark\n        $r.InsertAfter(\"DONE\")\n        break\n    }\n}\n\n# 2. Add a trailing space run at the end of the \"... users name\" paragraph.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -match \"users name\\s*\\r?$\") {\n        $r = $p.Range\n        $r.End = $r.End - 1   # exclude the paragraph mark\n        $r.InsertAfter(\" \")\n        break\n    }\n}\n\n# 3. Insert a blank paragraph and a new task paragraph after the\n#    \"Ask Rory about simulating PayPal Login\" paragraph.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Ask Rory about simulating PayPal Login*\") {\n        $r = $p.Range\n        $r.InsertParagraphAfter()\n        $blank = $d.Paragraphs.Item($i + 1)\n        $blank.Range.InsertParagraphAfter()\n        $textPara = $d.Paragraphs.Item($i + 2)\n        $textPara.Range.InsertBefore(\"Add info button to reappear instructions \")\n        break\n    }\n}\n"}
